# Add a new "MLS" row to the SPIS Use Case sheet (Peppol Code Lists v9.4)
# Columns: A=Use Case ID, B=Initial release, C=State, D=Deprecation release,
#          E=Removal Date, F=Comment

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "MLS"

# "9.4" must be stored as text (quote-prefixed), matching the original
# author typing '9.4 into a text-formatted cell.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "'9.4"

$ws.Range("C2").Value = "active"
$ws.Range("F2").Value = "TICC-410"

# Move the selection, as recorded in the saved workbook view state.
[void]$ws.Range("F3").Select()
